$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update query text in B2, B3, B4 to add ORDER BY / LIMIT clauses ---

$b2 = $ws.Range("B2").Value()
$b2 = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100"
$ws.Range("B2").Value = $b2

$b3 = $ws.Range("B3").Value()
$b3 = $b3 + "`n order By samp.sample_id ASC LIMIT 100"
$ws.Range("B3").Value = $b3

$b4 = $ws.Range("B4").Value()
$b4 = $b4.Replace("    order by f.file_name", "     order By f.file_name ASC LIMIT 100")
$ws.Range("B4").Value = $b4

# --- Update selected cell from D4 to B4 ---
$ws.Range("B4").Select()

# --- Row heights grow to fit the newly-added wrapped text (Excel auto-fit) ---
$ws.Rows.Item(2).RowHeight = 360
$ws.Rows.Item(3).RowHeight = 374.4
